$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (Cena shifts to D, Powierzchnia shifts to E)
$ws.Columns("C").Insert()

# New header cell C1 = "Data"
$ws.Range("C1").Value = "Data"

# Formulas for the new column referencing column A instead of C
$ws.Range("C2").Formula = "=INDIRECT(`"'`"&`$A2&`"'!A2`")"
$ws.Range("C3").Formula = "=INDIRECT(`"'`"&`$A3&`"'!A2`")"
